# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Changes to the "Hoja1" Estado de Cuenta worksheet:
#  1. Remove the worker "ALONSO FLOREZ CARPIO" (doc 3821319) entirely
#     -- this was row 23, a single period (2412) entry -- by deleting
#     that whole row (shifts everything below up by one, which also
#     naturally pulls the signature-block rows 35/36 up to 34/35).
#  2. Re-sort / rewrite the remaining detail rows so periods run in
#     ascending order (2412..2506) with each period showing first the
#     JHON EIDER CASTAÑEDA SAMUDIO row, then the IVETH QUIÑONES GIL row.
#  3. Update the VALOR MORA total (E11) and Cant. Trabajadores (C13).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# 1. Drop the ALONSO FLOREZ CARPIO record (row 23) completely.
$ws.Rows.Item(23).Delete()

# 2. Update summary fields.
$ws.Range("E11").Value = 945952
$ws.Range("C13").Value = 2

# 3. Rewrite the detail table (rows 16-29) in the final order.
$data = @(
    @("CC", "1052216389", "JHON EIDER CASTAÑEDA SAMUDIO", "2412", 52000,  1300000),
    @("CC", "26794952",   "IVETH QUIÑONES GIL",           "2412", 85760,  2144000),
    @("CC", "1052216389", "JHON EIDER CASTAÑEDA SAMUDIO", "2501", 52000,  1300000),
    @("CC", "26794952",   "IVETH QUIÑONES GIL",           "2501", 85760,  2144000),
    @("CC", "1052216389", "JHON EIDER CASTAÑEDA SAMUDIO", "2502", 52000,  1300000),
    @("CC", "26794952",   "IVETH QUIÑONES GIL",           "2502", 85760,  2144000),
    @("CC", "1052216389", "JHON EIDER CASTAÑEDA SAMUDIO", "2503", 52000,  1300000),
    @("CC", "26794952",   "IVETH QUIÑONES GIL",           "2503", 85760,  2144000),
    @("CC", "1052216389", "JHON EIDER CASTAÑEDA SAMUDIO", "2504", 52000,  1300000),
    @("CC", "26794952",   "IVETH QUIÑONES GIL",           "2504", 85760,  2144000),
    @("CC", "1052216389", "JHON EIDER CASTAÑEDA SAMUDIO", "2505", 52000,  1300000),
    @("CC", "26794952",   "IVETH QUIÑONES GIL",           "2505", 85760,  2144000),
    @("CC", "1052216389", "JHON EIDER CASTAÑEDA SAMUDIO", "2506", 45066,  1300000),
    @("CC", "26794952",   "IVETH QUIÑONES GIL",           "2506", 74326,  2144000)
)

$startRow = 16
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]
    $ws.Range("B$row").Value = $rec[0]
    $ws.Range("C$row").Value = $rec[1]
    $ws.Range("D$row").Value = $rec[2]
    $ws.Range("E$row").Value = $rec[3]
    $ws.Range("F$row").Value = $rec[4]
    $ws.Range("G$row").Value = $rec[5]
}
